$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Styling: the bold / bordered / centered header style (style index 1)
# is no longer used anywhere in the sheet, so strip it back to Normal and
# blank out A1's stray "Unnamed: 0" label (pandas artifact).
$ws.Range("A1").ClearContents()
$ws.Range("A1:AJ1").Style = "Normal"

# --- Corrected data-cleaning values (pre/post/total fixation metrics) ---

# Row 3 - Revisit count
$ws.Range("D3").Value = 5
$ws.Range("I3").Value = 10
$ws.Range("K3").Value = 24
$ws.Range("V3").Value = 25
$ws.Range("AF3").Value = 30

# Row 4 - Fixation count
$ws.Range("D4").Value = 8
$ws.Range("I4").Value = 24
$ws.Range("K4").Value = 208
$ws.Range("V4").Value = 153
$ws.Range("AF4").Value = 108

# Row 5 - Dwell time (ms)
$ws.Range("D5").Value = 3336.11
$ws.Range("I5").Value = 7107.01
$ws.Range("K5").Value = 50271.53
$ws.Range("V5").Value = 37641.72
$ws.Range("AF5").Value = 38789.16

# Row 6 - Dwell time (%)
$ws.Range("B6").Value = 0.13
$ws.Range("D6").Value = 2.53
$ws.Range("E6").Value = 1.98
$ws.Range("F6").Value = 1.27
$ws.Range("G6").Value = 0.76
$ws.Range("H6").Value = 0.42
$ws.Range("I6").Value = 5.4
$ws.Range("J6").Value = 8.75
$ws.Range("K6").Value = 38.17
$ws.Range("L6").Value = 0.76
$ws.Range("M6").Value = 13.73
$ws.Range("N6").Value = 5.73
$ws.Range("O6").Value = 3.9
$ws.Range("P6").Value = 8.029999999999999
$ws.Range("Q6").Value = 1.63
$ws.Range("R6").Value = 0.14
$ws.Range("U6").Value = 0.62
$ws.Range("V6").Value = 28.58
$ws.Range("W6").Value = 4.66
$ws.Range("X6").Value = 0.11
$ws.Range("Y6").Value = 2.7
$ws.Range("Z6").Value = 0.37
$ws.Range("AA6").Value = 3.43
$ws.Range("AB6").Value = 1.08
$ws.Range("AC6").Value = 3.29
$ws.Range("AD6").Value = 7.11
$ws.Range("AE6").Value = 0.25
$ws.Range("AF6").Value = 29.45
$ws.Range("AG6").Value = 0.06
$ws.Range("AH6").Value = 0.34
$ws.Range("AI6").Value = 1.22
$ws.Range("AJ6").Value = 3.55

# Row 7 - Fixation duration (ms)
$ws.Range("D7").Value = 417.01
$ws.Range("I7").Value = 296.13
$ws.Range("K7").Value = 241.69
$ws.Range("V7").Value = 246.02
$ws.Range("AF7").Value = 359.16

# --- Row 10 was a stray duplicate-of-blank row; drop it entirely ---
$ws.Rows.Item(10).Delete()
